$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "Quantity"
$ws.Range("A1:B1").Font.Bold = $true

# File rows
$files = @(
    "level1.DXF",
    "level2.DXF",
    "level3.DXF",
    "level4.DXF",
    "level5.DXF",
    "level6.DXF",
    "control_panel.DXF",
    "LIDAR_plate.DXF"
)

$row = 2
foreach ($f in $files) {
    $ws.Cells.Item($row, 1).Value = $f
    $ws.Cells.Item($row, 2).Value = 1
    $row++
}

# Autofit column A to match bestFit width
$ws.Columns.Item(1).AutoFit() | Out-Null

# Selection matches diff's saved selection state
$ws.Range("B10").Select()

# Page setup orientation
$ws.PageSetup.Orientation = 1
